# Add newly-found source data for the 2010-2019 decade as a new column
# in the GHG emissions table, inserted between the "2000-2009" and
# "2012-2021" columns (i.e. before the current column F).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F; this shifts the existing
# "2012-2021", "2021" and "2022 (projection)" columns one place to the
# right (F->G, G->H, H->I) and carries their formatting/values along.
$ws.Columns("F:F").Insert()

# Populate the freshly inserted column F with the 2010-2019 data.
$ws.Range("F1").Value = "2010-2019"
$ws.Range("F2").Value = "53±5.3"
$ws.Range("F3").Value = "36±2.8"
$ws.Range("F4").Value = "4.7±3.3"
$ws.Range("F5").Value = "8.6±2.6"
$ws.Range("F6").Value = "2.7±1.6"
$ws.Range("F7").Value = "1.5±0.46"
